$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extra Wishlist")

# Row 21: Sony FE 40mm f/2.5 G
$ws.Range("A21").Value = "Sony FE 40mm f/2.5 G"
$ws.Range("B21").Value = "https://static01.galaxus.com/productimages/1/9/9/3/2/3/5/8/4/3/6/8/0/4/6/0/0/8/2/01981e19-c81d-77c1-a3c2-f578d3601927_2880.avif"
$ws.Range("C21").Value = "https://www.digitec.ch/en/s1/product/sony-fe-40mm-f25-g-sony-e-full-size-lenses-15382655"
$ws.Range("D21").Value = 520

# Row 22: Sony FE 35mm f/1.4 GM
$ws.Range("A22").Value = "Sony FE 35mm f/1.4 GM"
$ws.Range("B22").Value = "https://static01.galaxus.com/productimages/4/1/4/3/9/0/3/8/1%20Digitec_A-Mid%203.jpg_2880.avif"
$ws.Range("C22").Value = "https://www.digitec.ch/en/s1/product/sony-fe-35mm-f14-gm-sony-e-full-size-lenses-14596684"
$ws.Range("D22").Value = 1236

$ws.Range("D23").Select()
